$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on cells whose new values look numeric,
# so Excel keeps them as text (matching the original inline-string cells)
# instead of auto-converting them to numbers.
$textCells = @("D5", "D8", "D9", "D10", "D11", "D12", "D14", "D17", "D19", "D21", "D22", "D25", "D26", "D27", "D28", "D29", "D33", "D34", "D35", "D36", "D37", "D40", "D41", "D42", "D43", "D44", "D46", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '35.232.29'
$ws.Range('E2').Value = '  -0.89%  '
$ws.Range('D3').Value = '1.896.35'
$ws.Range('E3').Value = '  -0.40%  '
$ws.Range('E4').Value = '  -0.20%  '
$ws.Range('D5').Value = '246.17'
$ws.Range('E5').Value = '  -0.07%  '
$ws.Range('E6').Value = '  +9.16%  '
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('D8').Value = '40.46'
$ws.Range('E8').Value = '  -4.21%  '
$ws.Range('D9').Value = '0.347'
$ws.Range('E9').Value = '  +2.39%  '
$ws.Range('D10').Value = '52.33'
$ws.Range('E10').Value = '  +8.56%  '
$ws.Range('D11').Value = '0.0719'
$ws.Range('E11').Value = '  +2.12%  '
$ws.Range('D12').Value = '0.0987'
$ws.Range('E12').Value = '  -1.17%  '
$ws.Range('D13').Value = '2.172.65'
$ws.Range('D14').Value = '12.55'
$ws.Range('E14').Value = '  +1.03%  '
$ws.Range('E15').Value = '  +2.30%  '
$ws.Range('D16').Value = '1.874.42'
$ws.Range('E16').Value = '  -2.03%  '
$ws.Range('D17').Value = '4.79'
$ws.Range('E17').Value = '  -1.62%  '
$ws.Range('D18').Value = '35.236.48'
$ws.Range('E18').Value = '  -0.82%  '
$ws.Range('D19').Value = '71.95'
$ws.Range('E19').Value = '  -0.01%  '
$ws.Range('D20').Value = '0.0₃0817'
$ws.Range('E20').Value = '  +0.67%  '
$ws.Range('D21').Value = '240.45'
$ws.Range('E21').Value = '  -1.35%  '
$ws.Range('D22').Value = '12.68'
$ws.Range('E22').Value = '  +1.65%  '
$ws.Range('E23').Value = '  -2.17%  '
$ws.Range('E24').Value = '  -0.23%  '
$ws.Range('D25').Value = '2.33'
$ws.Range('E25').Value = '  +1.87%  '
$ws.Range('D26').Value = '2.35'
$ws.Range('E26').Value = '  +5.08%  '
$ws.Range('D27').Value = '168.11'
$ws.Range('E27').Value = '  -2.36%  '
$ws.Range('D28').Value = '8.63'
$ws.Range('E28').Value = '  +0.50%  '
$ws.Range('D29').Value = '18.35'
$ws.Range('E29').Value = '  +1.97%  '
$ws.Range('E30').Value = '  +4.31%  '
$ws.Range('E31').Value = '  +20.48%  '
$ws.Range('E32').Value = '  +1.16%  '
$ws.Range('D33').Value = '0.0566'
$ws.Range('E33').Value = '  -0.03%  '
$ws.Range('B34').Value = 'BinanceUSD'
$ws.Range('C34').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D34').Value = '1.01'
$ws.Range('E34').Value = '  -0.09%  '
$ws.Range('B35').Value = 'WEMIXToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D35').Value = '1.87'
$ws.Range('E35').Value = '  +6.88%  '
$ws.Range('D36').Value = '4.09'
$ws.Range('E36').Value = '  -1.77%  '
$ws.Range('D37').Value = '1.51'
$ws.Range('E37').Value = '  +15.38%  '
$ws.Range('E38').Value = '  -8.53%  '
$ws.Range('E39').Value = '  -0.40%  '
$ws.Range('D40').Value = '0.0655'
$ws.Range('E40').Value = '  +10.41%  '
$ws.Range('B41').Value = 'ARBITRUM'
$ws.Range('C41').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D41').Value = '1.09'
$ws.Range('E41').Value = '  -1.53%  '
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').Value = '0.0207'
$ws.Range('E42').Value = '  +0.55%  '
$ws.Range('B43').Value = 'InjectiveProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D43').Value = '16.29'
$ws.Range('E43').Value = '  +5.35%  '
$ws.Range('D44').Value = '91.96'
$ws.Range('E44').Value = '  +0.98%  '
$ws.Range('D45').Value = '1.346.57'
$ws.Range('E45').Value = '  -0.50%  '
$ws.Range('D46').Value = '2.41'
$ws.Range('E46').Value = '  +2.97%  '
$ws.Range('E47').Value = '  +0.14%  '
$ws.Range('D48').Value = '2.78'
$ws.Range('E48').Value = '  +0.90%  '
$ws.Range('D49').Value = '45.37'
$ws.Range('E49').Value = '  -11.37%  '
$ws.Range('D50').Value = '12.10'
$ws.Range('E50').Value = '  -4.03%  '
$ws.Range('D51').Value = '6.45'
$ws.Range('E51').Value = '  -3.19%  '
